$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Timp3"
$ws.Cells.Item(2,3).Value = "Kdr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 86.43264233333333
$ws.Cells.Item(2,8).Value = 259.297927
$ws.Cells.Item(2,9).Value = 0.4989038832435519
$ws.Cells.Item(2,10).Value = 0.4989038832435519
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 127.6999736666667
$ws.Cells.Item(2,14).Value = 383.099921
$ws.Cells.Item(2,15).Value = 0.9554352891750322
$ws.Cells.Item(2,16).Value = 0.9554352891750322
$ws.Cells.Item(2,17).Value = 11037.44614990709
$ws.Cells.Item(2,18).Value = 99337.01534916378
$ws.Cells.Item(2,19).Value = 0.4766703759573495
$ws.Cells.Item(2,20).Value = 0.4766703759573496

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Timp3"
$ws.Cells.Item(3,3).Value = "Kdr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 86.43264233333333
$ws.Cells.Item(3,8).Value = 259.297927
$ws.Cells.Item(3,9).Value = 0.4989038832435519
$ws.Cells.Item(3,10).Value = 0.4989038832435519
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.4321196666666667
$ws.Cells.Item(3,14).Value = 1.296359
$ws.Cells.Item(3,15).Value = 0.003233065495828321
$ws.Cells.Item(3,16).Value = 0.003233065495828321
$ws.Cells.Item(3,17).Value = 37.34924459419922
$ws.Cells.Item(3,18).Value = 336.143201347793
$ws.Cells.Item(3,19).Value = 0.001612988930649489
$ws.Cells.Item(3,20).Value = 0.001612988930649489

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Timp3"
$ws.Cells.Item(4,3).Value = "Kdr"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 86.43264233333333
$ws.Cells.Item(4,8).Value = 259.297927
$ws.Cells.Item(4,9).Value = 0.4989038832435519
$ws.Cells.Item(4,10).Value = 0.4989038832435519
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 4.77305
$ws.Cells.Item(4,14).Value = 14.31915
$ws.Cells.Item(4,15).Value = 0.03571136528892854
$ws.Cells.Item(4,16).Value = 0.03571136528892854
$ws.Cells.Item(4,17).Value = 412.5473234891167
$ws.Cells.Item(4,18).Value = 3712.925911402051
$ws.Cells.Item(4,19).Value = 0.01781653881857543
$ws.Cells.Item(4,20).Value = 0.01781653881857544

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Timp3"
$ws.Cells.Item(5,3).Value = "Kdr"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 86.43264233333333
$ws.Cells.Item(5,8).Value = 259.297927
$ws.Cells.Item(5,9).Value = 0.4989038832435519
$ws.Cells.Item(5,10).Value = 0.4989038832435519
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.751186
$ws.Cells.Item(5,14).Value = 2.253558
$ws.Cells.Item(5,15).Value = 0.00562028004021099
$ws.Cells.Item(5,16).Value = 0.00562028004021099
$ws.Cells.Item(5,17).Value = 64.92699086380733
$ws.Cells.Item(5,18).Value = 584.342917774266
$ws.Cells.Item(5,19).Value = 0.002803979536977489
$ws.Cells.Item(5,20).Value = 0.002803979536977489

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Timp3"
$ws.Cells.Item(6,3).Value = "Kdr"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 31.60427266666667
$ws.Cells.Item(6,8).Value = 94.81281800000001
$ws.Cells.Item(6,9).Value = 0.1824252265675234
$ws.Cells.Item(6,10).Value = 0.1824252265675234
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 127.6999736666667
$ws.Cells.Item(6,14).Value = 383.099921
$ws.Cells.Item(6,15).Value = 0.9554352891750322
$ws.Cells.Item(6,16).Value = 0.9554352891750322
$ws.Cells.Item(6,17).Value = 4035.864787287487
$ws.Cells.Item(6,18).Value = 36322.78308558738
$ws.Cells.Item(6,19).Value = 0.1742954990983625
$ws.Cells.Item(6,20).Value = 0.1742954990983625

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Timp3"
$ws.Cells.Item(7,3).Value = "Kdr"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 31.60427266666667
$ws.Cells.Item(7,8).Value = 94.81281800000001
$ws.Cells.Item(7,9).Value = 0.1824252265675234
$ws.Cells.Item(7,10).Value = 0.1824252265675234
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.4321196666666667
$ws.Cells.Item(7,14).Value = 1.296359
$ws.Cells.Item(7,15).Value = 0.003233065495828321
$ws.Cells.Item(7,16).Value = 0.003233065495828321
$ws.Cells.Item(7,17).Value = 13.65682776996245
$ws.Cells.Item(7,18).Value = 122.911449929662
$ws.Cells.Item(7,19).Value = 0.0005897927055841238
$ws.Cells.Item(7,20).Value = 0.0005897927055841238

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Timp3"
$ws.Cells.Item(8,3).Value = "Kdr"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 31.60427266666667
$ws.Cells.Item(8,8).Value = 94.81281800000001
$ws.Cells.Item(8,9).Value = 0.1824252265675234
$ws.Cells.Item(8,10).Value = 0.1824252265675234
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 4.77305
$ws.Cells.Item(8,14).Value = 14.31915
$ws.Cells.Item(8,15).Value = 0.03571136528892854
$ws.Cells.Item(8,16).Value = 0.03571136528892854
$ws.Cells.Item(8,17).Value = 150.8487736516334
$ws.Cells.Item(8,18).Value = 1357.6389628647
$ws.Cells.Item(8,19).Value = 0.006514653903868378
$ws.Cells.Item(8,20).Value = 0.006514653903868378

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Timp3"
$ws.Cells.Item(9,3).Value = "Kdr"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 31.60427266666667
$ws.Cells.Item(9,8).Value = 94.81281800000001
$ws.Cells.Item(9,9).Value = 0.1824252265675234
$ws.Cells.Item(9,10).Value = 0.1824252265675234
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.751186
$ws.Cells.Item(9,14).Value = 2.253558
$ws.Cells.Item(9,15).Value = 0.00562028004021099
$ws.Cells.Item(9,16).Value = 0.00562028004021099
$ws.Cells.Item(9,17).Value = 23.74068716738267
$ws.Cells.Item(9,18).Value = 213.666184506444
$ws.Cells.Item(9,19).Value = 0.001025280859708419
$ws.Cells.Item(9,20).Value = 0.001025280859708419

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Timp3"
$ws.Cells.Item(10,3).Value = "Kdr"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.02528933333333333
$ws.Cells.Item(10,8).Value = 0.075868
$ws.Cells.Item(10,9).Value = 0.0001459743247925071
$ws.Cells.Item(10,10).Value = 0.0001459743247925071
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 127.6999736666667
$ws.Cells.Item(10,14).Value = 383.099921
$ws.Cells.Item(10,15).Value = 0.9554352891750322
$ws.Cells.Item(10,16).Value = 0.9554352891750322
$ws.Cells.Item(10,17).Value = 3.229447200714222
$ws.Cells.Item(10,18).Value = 29.065024806428
$ws.Cells.Item(10,19).Value = 0.0001394690212202591
$ws.Cells.Item(10,20).Value = 0.0001394690212202591

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Timp3"
$ws.Cells.Item(11,3).Value = "Kdr"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.02528933333333333
$ws.Cells.Item(11,8).Value = 0.075868
$ws.Cells.Item(11,9).Value = 0.0001459743247925071
$ws.Cells.Item(11,10).Value = 0.0001459743247925071
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.4321196666666667
$ws.Cells.Item(11,14).Value = 1.296359
$ws.Cells.Item(11,15).Value = 0.003233065495828321
$ws.Cells.Item(11,16).Value = 0.003233065495828321
$ws.Cells.Item(11,17).Value = 0.01092801829022222
$ws.Cells.Item(11,18).Value = 0.09835216461200001
$ws.Cells.Item(11,19).Value = 0.0000004719445527634913
$ws.Cells.Item(11,20).Value = 0.0000004719445527634913

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Timp3"
$ws.Cells.Item(12,3).Value = "Kdr"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.02528933333333333
$ws.Cells.Item(12,8).Value = 0.075868
$ws.Cells.Item(12,9).Value = 0.0001459743247925071
$ws.Cells.Item(12,10).Value = 0.0001459743247925071
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 4.77305
$ws.Cells.Item(12,14).Value = 14.31915
$ws.Cells.Item(12,15).Value = 0.03571136528892854
$ws.Cells.Item(12,16).Value = 0.03571136528892854
$ws.Cells.Item(12,17).Value = 0.1207072524666667
$ws.Cells.Item(12,18).Value = 1.0863652722
$ws.Cells.Item(12,19).Value = 0.000005212942435469918
$ws.Cells.Item(12,20).Value = 0.000005212942435469918

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Timp3"
$ws.Cells.Item(13,3).Value = "Kdr"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.02528933333333333
$ws.Cells.Item(13,8).Value = 0.075868
$ws.Cells.Item(13,9).Value = 0.0001459743247925071
$ws.Cells.Item(13,10).Value = 0.0001459743247925071
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.751186
$ws.Cells.Item(13,14).Value = 2.253558
$ws.Cells.Item(13,15).Value = 0.00562028004021099
$ws.Cells.Item(13,16).Value = 0.00562028004021099
$ws.Cells.Item(13,17).Value = 0.01899699314933333
$ws.Cells.Item(13,18).Value = 0.170972938344
$ws.Cells.Item(13,19).Value = 0.0000008204165840146038
$ws.Cells.Item(13,20).Value = 0.0000008204165840146038

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Timp3"
$ws.Cells.Item(14,3).Value = "Kdr"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 55.182874
$ws.Cells.Item(14,8).Value = 165.548622
$ws.Cells.Item(14,9).Value = 0.3185249158641322
$ws.Cells.Item(14,10).Value = 0.3185249158641322
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 127.6999736666667
$ws.Cells.Item(14,14).Value = 383.099921
$ws.Cells.Item(14,15).Value = 0.9554352891750322
$ws.Cells.Item(14,16).Value = 0.9554352891750322
$ws.Cells.Item(14,17).Value = 7046.851556650984
$ws.Cells.Item(14,18).Value = 63421.66400985886
$ws.Cells.Item(14,19).Value = 0.3043299450981
$ws.Cells.Item(14,20).Value = 0.3043299450981

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Timp3"
$ws.Cells.Item(15,3).Value = "Kdr"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 55.182874
$ws.Cells.Item(15,8).Value = 165.548622
$ws.Cells.Item(15,9).Value = 0.3185249158641322
$ws.Cells.Item(15,10).Value = 0.3185249158641322
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.4321196666666667
$ws.Cells.Item(15,14).Value = 1.296359
$ws.Cells.Item(15,15).Value = 0.003233065495828321
$ws.Cells.Item(15,16).Value = 0.003233065495828321
$ws.Cells.Item(15,17).Value = 23.84560511858867
$ws.Cells.Item(15,18).Value = 214.610446067298
$ws.Cells.Item(15,19).Value = 0.001029811915041945
$ws.Cells.Item(15,20).Value = 0.001029811915041945

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Timp3"
$ws.Cells.Item(16,3).Value = "Kdr"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 55.182874
$ws.Cells.Item(16,8).Value = 165.548622
$ws.Cells.Item(16,9).Value = 0.3185249158641322
$ws.Cells.Item(16,10).Value = 0.3185249158641322
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 4.77305
$ws.Cells.Item(16,14).Value = 14.31915
$ws.Cells.Item(16,15).Value = 0.03571136528892854
$ws.Cells.Item(16,16).Value = 0.03571136528892854
$ws.Cells.Item(16,17).Value = 263.3906167457
$ws.Cells.Item(16,18).Value = 2370.5155507113
$ws.Cells.Item(16,19).Value = 0.01137495962404925
$ws.Cells.Item(16,20).Value = 0.01137495962404925

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Timp3"
$ws.Cells.Item(17,3).Value = "Kdr"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 55.182874
$ws.Cells.Item(17,8).Value = 165.548622
$ws.Cells.Item(17,9).Value = 0.3185249158641322
$ws.Cells.Item(17,10).Value = 0.3185249158641322
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.751186
$ws.Cells.Item(17,14).Value = 2.253558
$ws.Cells.Item(17,15).Value = 0.00562028004021099
$ws.Cells.Item(17,16).Value = 0.00562028004021099
$ws.Cells.Item(17,17).Value = 41.452602388564
$ws.Cells.Item(17,18).Value = 373.073421497076
$ws.Cells.Item(17,19).Value = 0.001790199226941067
$ws.Cells.Item(17,20).Value = 0.001790199226941067
